# Auto-generated Excel COM-interop script applying the Cactuar_Profits diff.
# Updates 200 numeric cells across the 8 leveling-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to refresh currentAveragePrice / LevePrice / LeveProfit market-data snapshot columns (H:N).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 388.05554
$ws.Cells.Item(9, 9).Value = 215.58333
$ws.Cells.Item(9, 11).Value = 215.58333
$ws.Cells.Item(9, 13).Value = -46.58332999999999

$ws.Cells.Item(82, 8).Value = 6813.8184
$ws.Cells.Item(82, 9).Value = 4158.6665
$ws.Cells.Item(82, 11).Value = 12475.9995
$ws.Cells.Item(82, 13).Value = -12069.9995

$ws.Cells.Item(85, 8).Value = 6813.8184
$ws.Cells.Item(85, 9).Value = 4158.6665
$ws.Cells.Item(85, 11).Value = 12475.9995
$ws.Cells.Item(85, 13).Value = -11071.9995

$ws.Cells.Item(98, 8).Value = 1766.9584
$ws.Cells.Item(98, 9).Value = 1713.3478
$ws.Cells.Item(98, 11).Value = 1713.3478
$ws.Cells.Item(98, 13).Value = -215.3478

$ws.Cells.Item(122, 8).Value = 1766.9584
$ws.Cells.Item(122, 9).Value = 1713.3478
$ws.Cells.Item(122, 11).Value = 5140.0434
$ws.Cells.Item(122, 13).Value = -2690.0434

$ws.Cells.Item(129, 8).Value = 3377.8
$ws.Cells.Item(129, 9).Value = 798.2
$ws.Cells.Item(129, 11).Value = 2394.6
$ws.Cells.Item(129, 13).Value = 2605.4

$ws.Cells.Item(132, 8).Value = 170816.78
$ws.Cells.Item(132, 9).Value = 254328.47
$ws.Cells.Item(132, 11).Value = 762985.41
$ws.Cells.Item(132, 13).Value = -760455.41

$ws.Cells.Item(133, 8).Value = 100780
$ws.Cells.Item(133, 10).Value = 100780
$ws.Cells.Item(133, 12).Value = 100780
$ws.Cells.Item(133, 14).Value = -110900

$ws.Cells.Item(135, 8).Value = 2927.926
$ws.Cells.Item(135, 9).Value = 1610
$ws.Cells.Item(135, 11).Value = 14490
$ws.Cells.Item(135, 13).Value = -11955

$ws.Cells.Item(137, 8).Value = 14930332
$ws.Cells.Item(137, 9).Value = 558513.4
$ws.Cells.Item(137, 10).Value = 66668880
$ws.Cells.Item(137, 11).Value = 1675540.2
$ws.Cells.Item(137, 12).Value = 200006640
$ws.Cells.Item(137, 13).Value = -1672990.2
$ws.Cells.Item(137, 14).Value = -200011740

$ws.Cells.Item(141, 8).Value = 4094.3948
$ws.Cells.Item(141, 9).Value = 4212.839
$ws.Cells.Item(141, 11).Value = 12638.517
$ws.Cells.Item(141, 13).Value = -7458.517


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16475.225
$ws.Cells.Item(32, 9).Value = 17735.645
$ws.Cells.Item(32, 11).Value = 17735.645
$ws.Cells.Item(32, 13).Value = -17448.645

$ws.Cells.Item(61, 8).Value = 5114.35
$ws.Cells.Item(61, 9).Value = 5229
$ws.Cells.Item(61, 11).Value = 5229
$ws.Cells.Item(61, 13).Value = -5017

$ws.Cells.Item(97, 8).Value = 1982
$ws.Cells.Item(97, 9).Value = 1224.75
$ws.Cells.Item(97, 11).Value = 1224.75
$ws.Cells.Item(97, 13).Value = -728.75

$ws.Cells.Item(105, 8).Value = 91000
$ws.Cells.Item(105, 10).Value = 91000
$ws.Cells.Item(105, 12).Value = 91000
$ws.Cells.Item(105, 14).Value = -97988

$ws.Cells.Item(136, 8).Value = 5114.35
$ws.Cells.Item(136, 9).Value = 5229
$ws.Cells.Item(136, 11).Value = 15687
$ws.Cells.Item(136, 13).Value = -13137


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 308
$ws.Cells.Item(16, 9).Value = 308
$ws.Cells.Item(16, 11).Value = 308
$ws.Cells.Item(16, 13).Value = -138

$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Range("N23").ClearContents()

$ws.Cells.Item(86, 8).Value = 1076.8334
$ws.Cells.Item(86, 9).Value = 1076.8334
$ws.Cells.Item(86, 11).Value = 1076.8334
$ws.Cells.Item(86, 13).Value = 46.16660000000002

$ws.Cells.Item(89, 8).Value = 1076.8334
$ws.Cells.Item(89, 9).Value = 1076.8334
$ws.Cells.Item(89, 11).Value = 5384.166999999999
$ws.Cells.Item(89, 13).Value = 231.8330000000005

$ws.Cells.Item(94, 8).Value = 1777.5238
$ws.Cells.Item(94, 9).Value = 1441.5625
$ws.Cells.Item(94, 10).Value = 2852.6
$ws.Cells.Item(94, 11).Value = 1441.5625
$ws.Cells.Item(94, 12).Value = 2852.6
$ws.Cells.Item(94, 13).Value = -990.5625
$ws.Cells.Item(94, 14).Value = -3754.6

$ws.Cells.Item(107, 8).Value = 2352.6785
$ws.Cells.Item(107, 9).Value = 1294.409
$ws.Cells.Item(107, 11).Value = 1294.409
$ws.Cells.Item(107, 13).Value = 625.5909999999999

$ws.Cells.Item(134, 8).Value = 4749.375
$ws.Cells.Item(134, 9).Value = 4332.5
$ws.Cells.Item(134, 10).Value = 6000
$ws.Cells.Item(134, 11).Value = 12997.5
$ws.Cells.Item(134, 12).Value = 18000
$ws.Cells.Item(134, 13).Value = -10462.5
$ws.Cells.Item(134, 14).Value = -23070


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 17547966
$ws.Cells.Item(31, 9).Value = 43480824
$ws.Cells.Item(31, 11).Value = 43480824
$ws.Cells.Item(31, 13).Value = -43480529

$ws.Cells.Item(34, 8).Value = 17547966
$ws.Cells.Item(34, 9).Value = 43480824
$ws.Cells.Item(34, 11).Value = 43480824
$ws.Cells.Item(34, 13).Value = -43480622

$ws.Cells.Item(58, 8).Value = 402027.53
$ws.Cells.Item(58, 9).Value = 2233.4736
$ws.Cells.Item(58, 11).Value = 2233.4736
$ws.Cells.Item(58, 13).Value = -2030.4736

$ws.Cells.Item(110, 8).Value = 70700.5
$ws.Cells.Item(110, 10).Value = 70700.5
$ws.Cells.Item(110, 12).Value = 70700.5
$ws.Cells.Item(110, 14).Value = -78880.5

$ws.Cells.Item(132, 8).Value = 30312250
$ws.Cells.Item(132, 9).Value = 43019668
$ws.Cells.Item(132, 10).Value = 9941
$ws.Cells.Item(132, 11).Value = 129059004
$ws.Cells.Item(132, 12).Value = 29823
$ws.Cells.Item(132, 13).Value = -129056474
$ws.Cells.Item(132, 14).Value = -34883

$ws.Cells.Item(136, 8).Value = 402027.53
$ws.Cells.Item(136, 9).Value = 2233.4736
$ws.Cells.Item(136, 11).Value = 6700.4208
$ws.Cells.Item(136, 13).Value = -4150.4208

$ws.Cells.Item(141, 8).Value = 88375.05
$ws.Cells.Item(141, 10).Value = 90255.664
$ws.Cells.Item(141, 12).Value = 90255.664
$ws.Cells.Item(141, 14).Value = -100615.664


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 2500.1538
$ws.Cells.Item(2, 9).Value = 159.2
$ws.Cells.Item(2, 10).Value = 3963.25
$ws.Cells.Item(2, 11).Value = 955.1999999999999
$ws.Cells.Item(2, 12).Value = 23779.5
$ws.Cells.Item(2, 13).Value = -842.1999999999999
$ws.Cells.Item(2, 14).Value = -24005.5

$ws.Cells.Item(47, 8).Value = 5788
$ws.Cells.Item(47, 9).Value = 5780
$ws.Cells.Item(47, 11).Value = 17340
$ws.Cells.Item(47, 13).Value = -16909

$ws.Cells.Item(114, 8).Value = 1681.8948
$ws.Cells.Item(114, 10).Value = 2370.5
$ws.Cells.Item(114, 12).Value = 7111.5
$ws.Cells.Item(114, 14).Value = -13619.5

$ws.Cells.Item(140, 8).Value = 11265.083
$ws.Cells.Item(140, 9).Value = 3024.8667
$ws.Cells.Item(140, 11).Value = 9074.6001
$ws.Cells.Item(140, 13).Value = -3894.6001

$ws.Cells.Item(141, 8).Value = 11655
$ws.Cells.Item(141, 10).Value = 21249
$ws.Cells.Item(141, 12).Value = 63747
$ws.Cells.Item(141, 14).Value = -74107


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 528766
$ws.Cells.Item(132, 9).Value = 118157.28
$ws.Cells.Item(132, 11).Value = 354471.84
$ws.Cells.Item(132, 13).Value = -351941.84


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 17333.334
$ws.Cells.Item(19, 9).Value = 42000
$ws.Cells.Item(19, 10).Value = 5000
$ws.Cells.Item(19, 11).Value = 42000
$ws.Cells.Item(19, 12).Value = 5000
$ws.Cells.Item(19, 13).Value = -41830
$ws.Cells.Item(19, 14).Value = -5340

$ws.Cells.Item(103, 8).Value = 28999.5
$ws.Cells.Item(103, 10).Value = 28999.5
$ws.Cells.Item(103, 12).Value = 28999.5
$ws.Cells.Item(103, 14).Value = -31343.5

$ws.Cells.Item(132, 8).Value = 5414.391
$ws.Cells.Item(132, 10).Value = 4099.9
$ws.Cells.Item(132, 12).Value = 12299.7
$ws.Cells.Item(132, 14).Value = -17359.7

$ws.Cells.Item(135, 8).Value = 130000
$ws.Cells.Item(135, 10).Value = 130000
$ws.Cells.Item(135, 12).Value = 130000
$ws.Cells.Item(135, 14).Value = -140140


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 7500000.5
$ws.Cells.Item(5, 9).Value = 10000001
$ws.Cells.Item(5, 10).Value = 5000000
$ws.Cells.Item(5, 11).Value = 10000001
$ws.Cells.Item(5, 12).Value = 5000000
$ws.Cells.Item(5, 13).Value = -9999889
$ws.Cells.Item(5, 14).Value = -5000224

$ws.Cells.Item(107, 8).Value = 530.36365
$ws.Cells.Item(107, 9).Value = 483.25
$ws.Cells.Item(107, 11).Value = 1449.75
$ws.Cells.Item(107, 13).Value = 470.25

$ws.Cells.Item(122, 8).Value = 3467.74
$ws.Cells.Item(122, 9).Value = 2677.054
$ws.Cells.Item(122, 10).Value = 5718.154
$ws.Cells.Item(122, 11).Value = 8031.162
$ws.Cells.Item(122, 12).Value = 17154.462
$ws.Cells.Item(122, 13).Value = -5581.162
$ws.Cells.Item(122, 14).Value = -22054.462

$ws.Cells.Item(132, 8).Value = 7048.375
$ws.Cells.Item(132, 9).Value = 3427.6
$ws.Cells.Item(132, 11).Value = 10282.8
$ws.Cells.Item(132, 13).Value = -7752.799999999999

